$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; existing row 18 and below shift down to 19+
$ws.Rows.Item(18).Insert()

# Populate the new row 18 with data (matches row above/below for constant columns)
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).Value = 44901
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100103
$ws.Cells.Item(18, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(18, 9).Value = 100103003
$ws.Cells.Item(18, 10).Value = "Damasco"
$ws.Cells.Item(18, 11).Value = "Castle Brite"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 100
$ws.Cells.Item(18, 14).Value = 24000
$ws.Cells.Item(18, 15).Value = 24000
$ws.Cells.Item(18, 16).Value = 24000
$ws.Cells.Item(18, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 19).Value = 1333
$ws.Cells.Item(18, 20).Value = 18
